$wb = $excel.ActiveWorkbook

# Updated "想去人数" (interest count) values for column F in sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$expoUpdates = @{
    5  = 15654
    6  = 418
    9  = 15426
    10 = 53
    11 = 9023
    12 = 382
    15 = 90
    16 = 198
    18 = 198
    20 = 54
    24 = 61
    25 = 1114
    26 = 2
    29 = 86
    31 = 41
    32 = 413
    35 = 253
    39 = 5548
}
foreach ($row in $expoUpdates.Keys) {
    $wsExpo.Range("F$row").Value = $expoUpdates[$row]
}

# Same updates mirrored on the "全部类型" (All types) sheet, which lists the
# same events but with a couple of extra rows, shifting some row numbers.
$wsAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    5  = 15654
    6  = 418
    9  = 15426
    10 = 53
    11 = 9023
    12 = 382
    15 = 90
    16 = 198
    18 = 198
    20 = 54
    24 = 61
    25 = 1114
    26 = 2
    29 = 86
    31 = 41
    34 = 413
    37 = 253
    41 = 5548
}
foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
